$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 rows before row 34 (pushes old row34(blank)+note(35)+... down by 3)
$ws.Rows("34:36").Insert()

# Row 34: new meeting entry 11/28 /1:00
$ws.Range("B34").Value = "11/28 /1:00"
$ws.Range("C34").Value = "Google Hangout"
$ws.Range("D34").Value = "A"
$ws.Range("E34").Value = "A"
$ws.Range("F34").Value = "A"
$ws.Range("G34").Value = "A"
$ws.Range("I34").Value = "A"

# Row 35: new meeting entry 12/1 /4:15
$ws.Range("B35").Value = "12/1 /4:15"
$ws.Range("C35").Value = "Google Hangout"
$ws.Range("D35").Value = "A"
$ws.Range("E35").Value = "A"
$ws.Range("F35").Value = "A"
$ws.Range("G35").Value = "A"
$ws.Range("I35").Value = "A"

# Row 36 & 37 are spacer rows under the new entries - keep only the
# border formatting on D:G and I (clear A, B, C, H, J, K formatting that
# Insert() propagated down from row 33)
foreach ($r in 36, 37) {
    $ws.Range("A$r").ClearFormats()
    $ws.Range("B$r").ClearFormats()
    $ws.Range("C$r").ClearFormats()
    $ws.Range("H$r").ClearFormats()
    $ws.Range("J$r").ClearFormats()
    $ws.Range("K$r").ClearFormats()
}

$ws.Range("A32").Select()

Write-Host "Done"
for ($r = 30; $r -le 42; $r++) {
    $v1 = $ws.Range("A$r").Value2
    $v2 = $ws.Range("B$r").Value2
    $v3 = $ws.Range("C$r").Value2
    Write-Host "Row $r A=$v1 B=$v2 C=$v3"
}
